$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34:113 down to 35:114.
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new record's data.
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44868
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 100112026
$ws.Range("G34").Value = "Haba"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 120
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = 11000
$ws.Range("N34").Value = "`$/saco 25 kilos"
$ws.Range("O34").Value = "Región del Maule"
$ws.Range("P34").Value = 440
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
